# Sync attendance_reports: swap the order of the "Recorded By" values in
# column G from "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"
# wherever that exact combination appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$searchRange = $ws.Range("G1:G$lastRow")

# Use Find/FindNext so only matching cells are touched (xlValues lookup,
# xlWhole match) instead of rewriting every cell in the column.
$firstAddress = $null
$found = $searchRange.Find($oldValue, [System.Reflection.Missing]::Value, -4163, 1, 1, 1, $false)

while ($found -ne $null) {
    if ($firstAddress -eq $null) {
        $firstAddress = $found.Address()
    } elseif ($found.Address() -eq $firstAddress) {
        break
    }

    $found.Value2 = $newValue

    $found = $searchRange.FindNext($found)
}
